$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "MSDS, MSDA, MBA in AI BigData, 일반 MBA 넷 중 어느 석사를 해야하나요?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/msds-msda-mba-which-one/#utm_source=rss&utm_medium=rss&utm_campaign=msds-msda-mba-which-one"

$ws.Range("D28").Value = "강화학습 환경들"
$ws.Range("E28").Value = "https://ropiens.tistory.com/109"

$ws.Range("D32").Value = "클래스를 이용하여 데커레이터 만들기"
$ws.Range("E32").Value = "https://dodonam.tistory.com/316"

$ws.Range("D36").Value = "Introduction to Scene Text Detection and Recognition"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/320"

$ws.Range("D51").Value = "[python] .ipynb 파일을 .py 파일로 변환하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/1173"
